$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Technologies line: reorder a few of the comma-separated entries.
# ---------------------------------------------------------------------------
$techFound = $d.Content.Find.Execute(
    " Node.js, Spring Boot, Docker, Kubernetes, Kafka, RabbitMQ, AWS, Redis, GRPC, REST APIs, Prometheus, Grafana",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Docker, Kubernetes, Kafka, Node.js, Spring Boot, RabbitMQ, AWS, Redis, GRPC, REST APIs, Prometheus, Grafana",
    2)
Write-Host "Technologies line updated:" $techFound

# ---------------------------------------------------------------------------
# 2. "Developed Proofs of Concepts ..." bullet -> expanded sentence about
#    Kafka based async communication. Use Range.Text assignment (rather than
#    Find/Replace) so the xml:space="preserve" attribute on <w:t> survives.
# ---------------------------------------------------------------------------
$pocRange = $d.Content
$found = $pocRange.Find.Execute("Developed Proofs of Concepts and performed performance, load and feasibility testing.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $pocRange.Text = "Developed Proofs of Concept and performed performance, load and feasibility testing of Kafka based async communication leading to switch over from synchronous to asynchronous communication between services in the ecosystem."
}
Write-Host "POC bullet updated:" $found

# ---------------------------------------------------------------------------
# 3. Replace the "Set up reporting and metrics dashboards ..." bullet with
#    five new bullets describing additional Angel One responsibilities.
#    We build a small OOXML fragment (WordprocessingML wrapped in the
#    flat-opc pkg:package envelope) and splice it in via Range.InsertXML so
#    we retain full control over each paragraph's pPr/rPr (in particular the
#    <w:u w:val="none"/> marker carried on four of the five paragraphs).
# ---------------------------------------------------------------------------
function New-BulletParagraphXml([string]$text, [bool]$underlineNoneOnMark) {
    if ($underlineNoneOnMark) {
        $pPr = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="360"/><w:rPr><w:u w:val="none"/></w:rPr></w:pPr>'
    } else {
        $pPr = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="360"/></w:pPr>'
    }
    return '<w:p>' + $pPr + '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + $text + '</w:t></w:r></w:p>'
}

$dashRange = $d.Content
$dashFound = $dashRange.Find.Execute("Set up reporting and metrics dashboards using prometheus and grafana, along with alerting setups on the same.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Dashboards bullet located:" $dashFound

if ($dashFound) {
    $newParagraphsXml = ""
    $newParagraphsXml += New-BulletParagraphXml "Performed domain-research, designed, reviewed and peer-reviewed domain data models and DB design." $true
    $newParagraphsXml += New-BulletParagraphXml "Created and presented release milestones and execution plans to streamline product delivery." $true
    $newParagraphsXml += New-BulletParagraphXml "Collaborated with relevant business stakeholders, operations teams and other teams to clarify requirements and define technical tasks for execution by self and other engineers." $true
    $newParagraphsXml += New-BulletParagraphXml "Built monitoring system based on Prometheus and Grafana, and alerting systems using alertmanager." $false
    $newParagraphsXml += New-BulletParagraphXml "Led multiple focus groups of junior engineers to execute and deliver multiple technical milestones across components." $true

    $flatOpc = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $newParagraphsXml + '</w:body>' + `
        '</w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $dashRange.InsertXML($flatOpc)
}
